# "Feito do 11 ao 14" - mark exercises 11 through 14 as done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Checkboxes for exercises 11-14 live in B12:B15 -> set to TRUE
$ws.Range("B12:B15").Value = $true

# Move the view/selection down to the next block (row 16)
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("A16").Select()
